$d = $word.ActiveDocument

$find = " correlation wityh steps, duration, and "
$replace = " correlation with steps, duration, and "

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
